# Update "想去人数" (F column) values per commit 456a3b4 (gh-pages output regeneration)
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1270
$ws.Range("F5").Value = 5715
$ws.Range("F6").Value = 1806
$ws.Range("F7").Value = 1806
$ws.Range("F8").Value = 6375
$ws.Range("F9").Value = 142
$ws.Range("F10").Value = 1941
$ws.Range("F11").Value = 520
$ws.Range("F14").Value = 35
$ws.Range("F18").Value = 7975
$ws.Range("F19").Value = 7975
$ws.Range("F26").Value = 5
$ws.Range("F31").Value = 1786
$ws.Range("F32").Value = 808
$ws.Range("F33").Value = 382
$ws.Range("F36").Value = 9
$ws.Range("F37").Value = 85
$ws.Range("F39").Value = 3923

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 369
$ws.Range("F14").Value = 10
$ws.Range("F16").Value = 28

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 2283

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 2283
$ws.Range("F5").Value = 1270
$ws.Range("F9").Value = 369
$ws.Range("F10").Value = 5715
$ws.Range("F12").Value = 1806
$ws.Range("F13").Value = 1806
$ws.Range("F14").Value = 6375
$ws.Range("F15").Value = 142
$ws.Range("F16").Value = 1941
$ws.Range("F18").Value = 520
$ws.Range("F21").Value = 35
$ws.Range("F22").Value = 10
$ws.Range("F24").Value = 7975
$ws.Range("F25").Value = 7975
$ws.Range("F32").Value = 5
$ws.Range("F36").Value = 1786
$ws.Range("F37").Value = 808
$ws.Range("F39").Value = 382
$ws.Range("F41").Value = 28
$ws.Range("F47").Value = 3923

